$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper: remove bold formatting from a run of text found by exact
# search, without leaving a stray empty <w:rPr/> behind. We do this by
# deleting the found text and re-inserting it immediately before the
# following text (which merges it into the neighbouring, unbolded run).
# ------------------------------------------------------------------
function Remove-BoldExact($searchText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Text = ""
    $found = $rng.Find.Execute($searchText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Delete()
        $rng.InsertBefore($searchText)
    }
}

# 1) Remove bold from "twenty page plus documents"
Remove-BoldExact("twenty page plus documents")

# 2) "STILL" -> "still"
$d.Content.Find.Execute("STILL", $true, $false, $false, $false, $false, $true, 1, $false, "still", 2) | Out-Null

# 3) Remove bold from "you" (the one right before " need!")
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("you", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $after = $d.Range($rng.End, [Math]::Min($rng.End + 8, $d.Content.End))
    if ($after.Text -eq " need!  ") {
        $rng.Delete()
        $rng.InsertBefore("you")
        break
    }
    $rng.Collapse(0)
    $found = $rng.Find.Execute("you", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
}

# 4) Move the _GoBack bookmark to sit between "n" and "eed!" in " need!  "
$rng = $d.Content
$found = $rng.Find.Execute("need!")
if ($found) {
    $bmRange = $d.Range($rng.Start + 1, $rng.Start + 1)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}

# 5) Remove bold from "not"
Remove-BoldExact("not")

# 6) "EVERYTHING" -> "everything"
$d.Content.Find.Execute("EVERYTHING", $true, $false, $false, $false, $false, $true, 1, $false, "everything", 2) | Out-Null

$word.ActiveDocument.Saved = $false
